# Split the two long "Programa" run texts (PT + EN) into one <w:t> per
# numbered item, joined by manual line breaks (<w:br/>), instead of one
# giant run of concatenated sentences. We insert a "^l" (manual line
# break) immediately before each "N)" item marker (2..12); Word's find
# engine splits the owning run's text node at that point and inserts a
# <w:br/>, which is exactly the XML shape the diff shows, and formatting
# (e.g. the italic rPr on the English paragraph) is preserved on the run.
$d = $word.ActiveDocument

# --- Portuguese paragraph (Programa) ---
$d.Content.Find.Execute("2)Amostragem: Amostr", $true, $false, $false, $false, $false, $true, 1, $false, "^l2)Amostragem: Amostr", 2) | Out-Null
$d.Content.Find.Execute("3)Conceitos de Proba", $true, $false, $false, $false, $false, $true, 1, $false, "^l3)Conceitos de Proba", 2) | Out-Null
$d.Content.Find.Execute("4)Variáveis Aleatóri", $true, $false, $false, $false, $false, $true, 1, $false, "^l4)Variáveis Aleatóri", 2) | Out-Null
$d.Content.Find.Execute("5)Variáveis Aleatóri", $true, $false, $false, $false, $false, $true, 1, $false, "^l5)Variáveis Aleatóri", 2) | Out-Null
$d.Content.Find.Execute("6)Aproximações: Apro", $true, $false, $false, $false, $false, $true, 1, $false, "^l6)Aproximações: Apro", 2) | Out-Null
$d.Content.Find.Execute("7)Teorema Central do", $true, $false, $false, $false, $false, $true, 1, $false, "^l7)Teorema Central do", 2) | Out-Null
$d.Content.Find.Execute("8)Conceitos de Teste", $true, $false, $false, $false, $false, $true, 1, $false, "^l8)Conceitos de Teste", 2) | Out-Null
$d.Content.Find.Execute("9)Testes de Hipótese", $true, $false, $false, $false, $false, $true, 1, $false, "^l9)Testes de Hipótese", 2) | Out-Null
$d.Content.Find.Execute("10)Testes de Hipótes", $true, $false, $false, $false, $false, $true, 1, $false, "^l10)Testes de Hipótes", 2) | Out-Null
$d.Content.Find.Execute("11) Análise de Variâ", $true, $false, $false, $false, $false, $true, 1, $false, "^l11) Análise de Variâ", 2) | Out-Null
$d.Content.Find.Execute("12)Regressão Linear ", $true, $false, $false, $false, $false, $true, 1, $false, "^l12)Regressão Linear ", 2) | Out-Null

# --- English paragraph (Programa, italic) ---
$d.Content.Find.Execute("2)Sampling methods: ", $true, $false, $false, $false, $false, $true, 1, $false, "^l2)Sampling methods: ", 2) | Out-Null
$d.Content.Find.Execute("3)Introduction to pr", $true, $false, $false, $false, $false, $true, 1, $false, "^l3)Introduction to pr", 2) | Out-Null
$d.Content.Find.Execute("4)Discrete Random Va", $true, $false, $false, $false, $false, $true, 1, $false, "^l4)Discrete Random Va", 2) | Out-Null
$d.Content.Find.Execute("5)Continuous Random ", $true, $false, $false, $false, $false, $true, 1, $false, "^l5)Continuous Random ", 2) | Out-Null
$d.Content.Find.Execute("6)Approximations: Ap", $true, $false, $false, $false, $false, $true, 1, $false, "^l6)Approximations: Ap", 2) | Out-Null
$d.Content.Find.Execute("7)Central Limit Theo", $true, $false, $false, $false, $false, $true, 1, $false, "^l7)Central Limit Theo", 2) | Out-Null
$d.Content.Find.Execute("8)Hypothesis test co", $true, $false, $false, $false, $false, $true, 1, $false, "^l8)Hypothesis test co", 2) | Out-Null
$d.Content.Find.Execute("9)Hypothesis test fo", $true, $false, $false, $false, $false, $true, 1, $false, "^l9)Hypothesis test fo", 2) | Out-Null
$d.Content.Find.Execute("10)Hypothesis test f", $true, $false, $false, $false, $false, $true, 1, $false, "^l10)Hypothesis test f", 2) | Out-Null
$d.Content.Find.Execute("11)Analysis of varia", $true, $false, $false, $false, $false, $true, 1, $false, "^l11)Analysis of varia", 2) | Out-Null
$d.Content.Find.Execute("12)Simple linear reg", $true, $false, $false, $false, $false, $true, 1, $false, "^l12)Simple linear reg", 2) | Out-Null
